# Add a new literature entry (row 36) to the 0-literature-list workbook.
# "belkin, niyogi" (2005) - towards a theoretical foundation of
# laplacian-based manifold methods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=authors, B=year, C=title, D=doc, E=keywords, F=notes
# Write E/F/C/D before A/B so new shared-string entries land in the same
# order as the target workbook (keywords, notes, title were newly added;
# authors/doc reuse already-existing shared strings).
$ws.Range("E36").Value = "laplace-beltrami"
$ws.Range("F36").Value = "theoretical foundation of using graph laplacian"
$ws.Range("C36").Value = "towards a theoretical foundation of laplacian-based manifold methods"
$ws.Range("D36").Value = "paper"
$ws.Range("A36").Value = "belkin, niyogi"
$ws.Range("B36").Value = 2005

# Match formatting of the row above (left/vcenter aligned, same style index)
# instead of leaving the COM-default style on the freshly written cells.
$ws.Range("A35:F35").Copy()
$ws.Range("A36:F36").PasteSpecial(-4122)

# Reflect the saved viewport/selection state (scrolled down to show the
# newly added row, with F39 as the active cell) as closely as this host
# allows.
$win = $excel.ActiveWindow
$win.ScrollRow = 31
$win.ScrollColumn = 1
[void]$ws.Range("F39").Select()
